$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# generator sheet: extend RES model - P/Q columns + O7 value + selection
# ---------------------------------------------------------------------------
$gen = $wb.Worksheets.Item("generator")
$gen.Activate()

# P column (curtailment upper bound?) updates on rows 2-5 (stays 80 elsewhere)
$gen.Range("P2").Value = 90
$gen.Range("P3").Value = 90
$gen.Range("P4").Value = 90
$gen.Range("P5").Value = 90

# O7 changes from 0 to 10000
$gen.Range("O7").Value = 10000

# Q column: previously empty cells (style s=11) now get values with a new
# centered style (fontId/fillId/borderId unchanged, alignment center/center
# added - this creates a brand-new cellXfs entry).
$gen.Range("Q2").HorizontalAlignment = -4108   # xlCenter
$gen.Range("Q2").VerticalAlignment = -4108     # xlCenter
$gen.Range("Q2").Value = 200
$gen.Range("Q3").Value = 200
$gen.Range("Q4").Value = 200
$gen.Range("Q5").Value = 200
$gen.Range("Q6").Value = 0
$gen.Range("Q7").Value = 0
$gen.Range("Q8").Value = 0
$gen.Range("Q9").Value = 0
$gen.Range("Q10").Value = 0
$gen.Range("Q11").Value = 0
$gen.Range("Q12").Value = 0
$gen.Range("Q13").Value = 0

# Broadcast the Q2 formatting (fill/border preserved, alignment added) to the
# rest of the Q column via a format-only paste so no extra/orphan style gets
# interned in cellXfs.
$gen.Range("Q2").Copy()
$gen.Range("Q3:Q13").PasteSpecial(-4122)       # xlPasteFormats
$excel.CutCopyMode = $false

# Selection moves from I1 to B2:B3
$gen.Range("B2:B3").Select()

# ---------------------------------------------------------------------------
# demand sheet: becomes the selected/active tab, selection moves to C8
# ---------------------------------------------------------------------------
$demand = $wb.Worksheets.Item("demand")
$demand.Activate()
$demand.Range("C8").Select()
